$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A429").Value = "baby"
$ws.Range("B429").Value = "EA 23.114"
$ws.Range("D429").Value = "*r"
$ws.Range("E429").Value = "*r"
$ws.Range("F429").Value = "*r"
$ws.Range("J429").Value = "Bebê"
$ws.Range("K429").Value = "baby"
$ws.Range("L429").Value = "赤ちゃん"

$ws.Range("A430").Value = "shark"
$ws.Range("B430").Value = "EA 23.119"
$ws.Range("D430").Value = "Tubarão  "
$ws.Range("E430").Value = "shark"
$ws.Range("F430").Value = "鮫"

$ws.Range("A431").Value = "fish_army"
$ws.Range("B431").Value = "EA 23.119"
$ws.Range("D431").Value = "Peixinho  "
$ws.Range("E431").Value = "small fish"
$ws.Range("F431").Value = "小魚"

$ws.Range("A432").Value = "fish_kuma"
$ws.Range("B432").Value = "EA 23.119"
$ws.Range("D432").Value = "Peixe-Coroa  "
$ws.Range("E432").Value = "crownfish"
$ws.Range("F432").Value = "クマノミ"

$ws.Range("A433").Value = "fish_coral"
$ws.Range("B433").Value = "EA 23.119"
$ws.Range("D433").Value = "Peixe-Coral  "
$ws.Range("E433").Value = "coralfish"
$ws.Range("F433").Value = "コーラルフィッシュ"

$ws.Range("A434").Value = "fish_xida"
$ws.Range("B434").Value = "EA 23.119"
$ws.Range("D434").Value = "Xidazoon  "
$ws.Range("E434").Value = "xidazoon"
$ws.Range("F434").Value = "シダズーン"

$ws.Range("A435").Value = "fish_angel"
$ws.Range("B435").Value = "EA 23.119"
$ws.Range("D435").Value = "Peixe-Anjo  "
$ws.Range("E435").Value = "angelfish"
$ws.Range("F435").Value = "エンゼルフィッシュ"

$ws.Range("A436").Value = "fish_seabass"
$ws.Range("B436").Value = "EA 23.119"
$ws.Range("D436").Value = "Robalo  "
$ws.Range("E436").Value = "sea bass"
$ws.Range("F436").Value = "シーバス"

$ws.Range("A437").Value = "fish_ piranha"
$ws.Range("B437").Value = "EA 23.119"
$ws.Range("D437").Value = "Piranha  "
$ws.Range("E437").Value = "piranha"
$ws.Range("F437").Value = "ピラニア"

$ws.Range("A438").Value = "balloonfish"
$ws.Range("B438").Value = "EA 23.116"
$ws.Range("D438").Value = "Peixe-Balão  "
$ws.Range("E438").Value = "balloonfish"
$ws.Range("F438").Value = "ハリセンボン"

$ws.Range("A439").Value = "tako"
$ws.Range("B439").Value = "EA 23.117"
$ws.Range("D439").Value = "Polvo  "
$ws.Range("E439").Value = "octopus"
$ws.Range("F439").Value = "オクトパス"

$ws.Range("A440").Value = "ika"
$ws.Range("B440").Value = "EA 23.117"
$ws.Range("D440").Value = "Bebê Kraken  "
$ws.Range("E440").Value = "baby kraken"
$ws.Range("F440").Value = "クラーケンの赤子"

$ws.Range("A441").Value = "jellyfish"
$ws.Range("B441").Value = "EA 23.116"
$ws.Range("D441").Value = "Água-Viva  "
$ws.Range("E441").Value = "jellyfish"
$ws.Range("F441").Value = "クラゲ"

$ws.Range("A442").Value = "tacchan"
$ws.Range("B442").Value = "EA 23.115"
$ws.Range("D442").Value = "Cavalo-Marinho  "
$ws.Range("E442").Value = "seahorse"
$ws.Range("F442").Value = "タッキー"

$ws.Range("A443").Value = "snail_sea"
$ws.Range("B443").Value = "EA 23.116"
$ws.Range("D443").Value = "Caramujo-Marinho  "
$ws.Range("E443").Value = "sea snail"
$ws.Range("F443").Value = "海かたつむり"

$ws.Range("A444").Value = "turtle"
$ws.Range("B444").Value = "EA 23.112"
$ws.Range("D444").Value = "Tartaruga  "
$ws.Range("E444").Value = "turtle"
$ws.Range("F444").Value = "亀"

$ws.Range("A445").Value = "sea_anemone"
$ws.Range("B445").Value = "EA 23.113"
$ws.Range("D445").Value = "Anêmona-do-Mar  "
$ws.Range("E445").Value = "sea anemone"
$ws.Range("F445").Value = "海のアネモネ"

$ws.Range("A446").Value = "dragon_leafy"
$ws.Range("B446").Value = "EA 23.114"
$ws.Range("D446").Value = "Dragão-Marinho Folhado  "
$ws.Range("E446").Value = "leafy sea #ele3 dragon"
$ws.Range("F446").Value = "リーフィーシー#ele3ドラゴン"

$ws.Range("A447").Value = "quickling"
$ws.Range("B447").Value = "EA 23.114"
$ws.Range("D447").Value = "Velozinho  "
$ws.Range("E447").Value = "quickling"
$ws.Range("F447").Value = "クイックリング"

$ws.Range("A448").Value = "quickling_archer"
$ws.Range("B448").Value = "EA 23.115"
$ws.Range("D448").Value = "Arqueiro Velozinho  "
$ws.Range("E448").Value = "quickling archer"
$ws.Range("F448").Value = "クイックリングの弓使い"

$ws.Range("A449").Value = "gomachan"
$ws.Range("B449").Value = "EA 23.112"
$ws.Range("D449").Value = "Foca"
$ws.Range("E449").Value = "Seal"
$ws.Range("F449").Value = "アザラシ"

$ws.Range("J429").Select()
